# Natmi following Dr Hou advice
# Updates the LR-pairs sheet (Inhba -> Bambi) with the corrected per-cluster-pair
# sender/receiver statistics, expanding the data from a 3-row sample to the full
# 3x3 (sender cluster x receiver cluster) cross-product (rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Inhba -> Bambi)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Inhba"
$ws.Cells.Item(2, 3).Value = "Bambi"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 3.675031333333333
$ws.Cells.Item(2, 8).Value = 11.025094
$ws.Cells.Item(2, 9).Value = 0.2032371147293133
$ws.Cells.Item(2, 10).Value = 0.2032371147293133
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.035948
$ws.Cells.Item(2, 14).Value = 6.107844
$ws.Cells.Item(2, 15).Value = 0.4919823674428878
$ws.Cells.Item(2, 16).Value = 0.4919823674428878
$ws.Cells.Item(2, 17).Value = 7.482172693037333
$ws.Cells.Item(2, 18).Value = 67.339554237336
$ws.Cells.Item(2, 19).Value = 0.09998907685678936
$ws.Cells.Item(2, 20).Value = 0.09998907685678936

# Row 3: ECs -> FAPs (Inhba -> Bambi)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Inhba"
$ws.Cells.Item(3, 3).Value = "Bambi"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 3.675031333333333
$ws.Cells.Item(3, 8).Value = 11.025094
$ws.Cells.Item(3, 9).Value = 0.2032371147293133
$ws.Cells.Item(3, 10).Value = 0.2032371147293133
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.335597666666667
$ws.Cells.Item(3, 14).Value = 4.006793
$ws.Cells.Item(3, 15).Value = 0.3227442459227168
$ws.Cells.Item(3, 16).Value = 0.3227442459227168
$ws.Cells.Item(3, 17).Value = 4.908363273726889
$ws.Cells.Item(3, 18).Value = 44.17526946354199
$ws.Cells.Item(3, 19).Value = 0.0655936093368209
$ws.Cells.Item(3, 20).Value = 0.0655936093368209

# Row 4: ECs -> sCs (Inhba -> Bambi)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Inhba"
$ws.Cells.Item(4, 3).Value = "Bambi"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 3.675031333333333
$ws.Cells.Item(4, 8).Value = 11.025094
$ws.Cells.Item(4, 9).Value = 0.2032371147293133
$ws.Cells.Item(4, 10).Value = 0.2032371147293133
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.7667083333333333
$ws.Cells.Item(4, 14).Value = 2.300125
$ws.Cells.Item(4, 15).Value = 0.1852733866343954
$ws.Cells.Item(4, 16).Value = 0.1852733866343954
$ws.Cells.Item(4, 17).Value = 2.817677148527777
$ws.Cells.Item(4, 18).Value = 25.35909433675
$ws.Cells.Item(4, 19).Value = 0.03765442853570304
$ws.Cells.Item(4, 20).Value = 0.03765442853570304

# Row 5: FAPs -> ECs (Inhba -> Bambi)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Inhba"
$ws.Cells.Item(5, 3).Value = "Bambi"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 10.108494
$ws.Cells.Item(5, 8).Value = 30.325482
$ws.Cells.Item(5, 9).Value = 0.5590213983169419
$ws.Cells.Item(5, 10).Value = 0.5590213983169419
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.035948
$ws.Cells.Item(5, 14).Value = 6.107844
$ws.Cells.Item(5, 15).Value = 0.4919823674428878
$ws.Cells.Item(5, 16).Value = 0.4919823674428878
$ws.Cells.Item(5, 17).Value = 20.580368142312
$ws.Cells.Item(5, 18).Value = 185.223313280808
$ws.Cells.Item(5, 19).Value = 0.2750286709952026
$ws.Cells.Item(5, 20).Value = 0.2750286709952026

# Row 6: FAPs -> FAPs (Inhba -> Bambi)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Inhba"
$ws.Cells.Item(6, 3).Value = "Bambi"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 10.108494
$ws.Cells.Item(6, 8).Value = 30.325482
$ws.Cells.Item(6, 9).Value = 0.5590213983169419
$ws.Cells.Item(6, 10).Value = 0.5590213983169419
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.335597666666667
$ws.Cells.Item(6, 14).Value = 4.006793
$ws.Cells.Item(6, 15).Value = 0.3227442459227168
$ws.Cells.Item(6, 16).Value = 0.3227442459227168
$ws.Cells.Item(6, 17).Value = 13.500880999914
$ws.Cells.Item(6, 18).Value = 121.507928999226
$ws.Cells.Item(6, 19).Value = 0.1804209396544641
$ws.Cells.Item(6, 20).Value = 0.1804209396544641

# Row 7: FAPs -> sCs (Inhba -> Bambi)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Inhba"
$ws.Cells.Item(7, 3).Value = "Bambi"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 10.108494
$ws.Cells.Item(7, 8).Value = 30.325482
$ws.Cells.Item(7, 9).Value = 0.5590213983169419
$ws.Cells.Item(7, 10).Value = 0.5590213983169419
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.7667083333333333
$ws.Cells.Item(7, 14).Value = 2.300125
$ws.Cells.Item(7, 15).Value = 0.1852733866343954
$ws.Cells.Item(7, 16).Value = 0.1852733866343954
$ws.Cells.Item(7, 17).Value = 7.750266587250001
$ws.Cells.Item(7, 18).Value = 69.75239928525001
$ws.Cells.Item(7, 19).Value = 0.1035717876672751
$ws.Cells.Item(7, 20).Value = 0.1035717876672751

# Row 8: sCs -> ECs (Inhba -> Bambi)
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Inhba"
$ws.Cells.Item(8, 3).Value = "Bambi"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.298956
$ws.Cells.Item(8, 8).Value = 12.896868
$ws.Cells.Item(8, 9).Value = 0.2377414869537448
$ws.Cells.Item(8, 10).Value = 0.2377414869537448
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.035948
$ws.Cells.Item(8, 14).Value = 6.107844
$ws.Cells.Item(8, 15).Value = 0.4919823674428878
$ws.Cells.Item(8, 16).Value = 0.4919823674428878
$ws.Cells.Item(8, 17).Value = 8.752450870288
$ws.Cells.Item(8, 18).Value = 78.77205783259201
$ws.Cells.Item(8, 19).Value = 0.1169646195908958
$ws.Cells.Item(8, 20).Value = 0.1169646195908958

# Row 9: sCs -> FAPs (Inhba -> Bambi)
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Inhba"
$ws.Cells.Item(9, 3).Value = "Bambi"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.298956
$ws.Cells.Item(9, 8).Value = 12.896868
$ws.Cells.Item(9, 9).Value = 0.2377414869537448
$ws.Cells.Item(9, 10).Value = 0.2377414869537448
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.335597666666667
$ws.Cells.Item(9, 14).Value = 4.006793
$ws.Cells.Item(9, 15).Value = 0.3227442459227168
$ws.Cells.Item(9, 16).Value = 0.3227442459227168
$ws.Cells.Item(9, 17).Value = 5.741675602702667
$ws.Cells.Item(9, 18).Value = 51.67508042432401
$ws.Cells.Item(9, 19).Value = 0.07672969693143178
$ws.Cells.Item(9, 20).Value = 0.07672969693143178

# Row 10: sCs -> sCs (Inhba -> Bambi)
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Inhba"
$ws.Cells.Item(10, 3).Value = "Bambi"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.298956
$ws.Cells.Item(10, 8).Value = 12.896868
$ws.Cells.Item(10, 9).Value = 0.2377414869537448
$ws.Cells.Item(10, 10).Value = 0.2377414869537448
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.7667083333333333
$ws.Cells.Item(10, 14).Value = 2.300125
$ws.Cells.Item(10, 15).Value = 0.1852733866343954
$ws.Cells.Item(10, 16).Value = 0.1852733866343954
$ws.Cells.Item(10, 17).Value = 3.296045389833334
$ws.Cells.Item(10, 18).Value = 29.6644085085
$ws.Cells.Item(10, 19).Value = 0.04404717043141723
$ws.Cells.Item(10, 20).Value = 0.04404717043141723

Write-Output "Updated rows 2-10 (A1:T10) with corrected Inhba -> Bambi LR-pair statistics."
